$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new worker (ARMANDO RAFAEL GARCIA ALVIS, doc 1048935962) is inserted right
# after the first worker (ALAN JOSE GARCIA ALVIS) in the table, pushing the
# remaining workers down by one row. ARMANDO previously occupied the last row
# of the table (20); he now occupies row 17, and the others shift down.

$ws.Range("C17").Value = "1048935962"
$ws.Range("D17").Value = "ARMANDO RAFAEL GARCIA ALVIS"
$ws.Range("F17").Value = 36341

$ws.Range("C18").Value = "73507432"
$ws.Range("D18").Value = "DUGLAS MARTELO ATENCIO"
$ws.Range("F18").Value = 9691

$ws.Range("C19").Value = "1007874342"
$ws.Range("D19").Value = "ROQUE JACINTO GARCIA ALVIS"
$ws.Range("F19").Value = 9691

$ws.Range("C20").Value = "73508131"
$ws.Range("D20").Value = "LUIS GERMAN MARTELO JIMENEZ"
$ws.Range("F20").Value = 36341
